# due date added when recon_accord is true
# Insert a new row into the "on_recon" sheet (sheet index 8) for the
# "message.orders.settlements.due_date" field, just above the existing
# "message.orders.settlements.status" row, and move the active
# selection from the "recon" sheet to "on_recon".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# Make room for the new row (everything from row 28 down shifts by one).
$ws.Rows.Item(28).Insert() | Out-Null

$ws.Cells.Item(28, 1).Value = "message.orders.settlements.due_date"
$ws.Cells.Item(28, 2).Value = $false
$ws.Cells.Item(28, 3).Value = "string"
$ws.Cells.Item(28, 4).Value = "any"

# Description first so the shared-string table fills in the same order
# as the authored workbook (description before the sample value).
$ws.Cells.Item(28, 6).Value = "due date of settlement in case recon_accord is true"

# Match the formatting already used for the description column elsewhere
# on this sheet (style used by F22:F26) by copying it onto the new cell.
$ws.Cells.Item(26, 6).Copy()
$ws.Cells.Item(28, 6).PasteSpecial(-4122)

# The sample value is a plain string that looks like a date ("2024-05-09"),
# so force text formatting before assigning it, otherwise Excel will
# silently convert it to a date serial number.
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "2024-05-09"

# The workbook was left with "on_recon" as the active sheet/tab.
$ws.Activate() | Out-Null
$ws.Range("D18").Select() | Out-Null

Write-Output "done"
